$wb = $excel.ActiveWorkbook

# "Weekly Quantity" sheet: remove the row for 2023-06-25 week (qty 40)
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows.Item(9).Delete()

# "Monthly Trend" sheet: the June 2023 total must drop by 40 (560 -> 520)
# since the removed weekly entry belonged to that month.
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Cells.Item(5, 2).Value = 520
